$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear order-data cells B2:H4 but keep them present with the sheet's
# default ("Normal") cell style, matching how the rest of the sheet's
# untouched cells are styled.
$rng = $ws.Range("B2:H4")
$rng.Value = $null
$rng.Style = "Normal"

# A3 and A4 also lose their values (only A2's order-counter of 2 stays).
$ws.Range("A3:A4").Value = $null
$ws.Range("A3:A4").Style = "Normal"

# The used range grew by two more (blank) rows: row 5 has a single
# (empty) cell in column A, row 6 has no cells at all.
$ws.Range("A5").Value = $null
$ws.Range("A5").Style = "Normal"

$ws.Rows.Item(6).OutlineLevel = 0
